$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Direct text assignments (values Excel will not mis-parse as numbers) ---
$ws.Range('D2').Value = '71.764.57'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '3.995.93'
$ws.Range('E3').Value = '  -0.93%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('E6').Value = '  +1.28%  '
$ws.Range('E7').Value = '  +9.71%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('E10').Value = '  -4.24%  '
$ws.Range('E11').Value = '  -5.76%  '
$ws.Range('E12').Value = '  +1.63%  '
$ws.Range('E13').Value = '  -3.08%  '
$ws.Range('D14').Value = '4.640.95'
$ws.Range('E14').Value = '  -0.61%  '
$ws.Range('D15').Value = '3.999.61'
$ws.Range('E15').Value = '  -2.02%  '
$ws.Range('E16').Value = '  -3.02%  '
$ws.Range('E17').Value = '  -4.92%  '
$ws.Range('E18').Value = '  -1.02%  '
$ws.Range('E19').Value = '  -3.24%  '
$ws.Range('D20').Value = '71.670.89'
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('E21').Value = '  -4.59%  '
$ws.Range('E22').Value = '  +2.58%  '
$ws.Range('E23').Value = '  -4.51%  '
$ws.Range('E24').Value = '  +2.52%  '
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('E26').Value = '  -8.15%  '
$ws.Range('E27').Value = '  -4.42%  '
$ws.Range('E28').Value = '  +0.63%  '
$ws.Range('E29').Value = '  -2.41%  '
$ws.Range('E30').Value = '  +23.12%  '
$ws.Range('E31').Value = '  -3.66%  '
$ws.Range('E32').Value = '  -2.47%  '
$ws.Range('E33').Value = '  -3.64%  '
$ws.Range('E34').Value = '  -1.53%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('E35').Value = '  -4.09%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('E36').Value = '  +5.22%  '
$ws.Range('E37').Value = '  -3.76%  '
$ws.Range('E38').Value = '  -1.34%  '
$ws.Range('D39').Value = '0.0₃0822'
$ws.Range('E39').Value = '  -9.89%  '
$ws.Range('E40').Value = '  -3.52%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('E43').Value = '  -2.53%  '
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('E45').Value = '  +1.61%  '
$ws.Range('E46').Value = '  -3.99%  '
$ws.Range('E47').Value = '  -8.80%  '
$ws.Range('E48').Value = '  +2.04%  '
$ws.Range('E49').Value = '  -6.36%  '
$ws.Range('E50').Value = '  -3.22%  '
$ws.Range('E51').Value = '  +1.32%  '

# --- Numeric-looking D-column text values: stage in column ZZ as Text, then ---
# --- Copy + PasteSpecial (values+formats) into destination so Excel keeps them ---
# --- as text instead of auto-converting to numbers; then clear the staging cells. ---
$ws.Range('ZZ5:ZZ7').NumberFormat = "@"
$ws.Range('ZZ5').Value = '528.44'
$ws.Range('ZZ6').Value = '150.31'
$ws.Range('ZZ7').Value = '0.691'
$ws.Range('ZZ5:ZZ7').Copy()
$ws.Range('D5:D7').PasteSpecial(-4104, $null, $null, $null)
$ws.Range('ZZ5:ZZ7').Clear()

$ws.Range('ZZ9:ZZ13').NumberFormat = "@"
$ws.Range('ZZ9').Value = '0.742'
$ws.Range('ZZ10').Value = '0.170'
$ws.Range('ZZ11').Value = '0.0000326'
$ws.Range('ZZ12').Value = '47.36'
$ws.Range('ZZ13').Value = '10.59'
$ws.Range('ZZ9:ZZ13').Copy()
$ws.Range('D9:D13').PasteSpecial(-4104, $null, $null, $null)
$ws.Range('ZZ9:ZZ13').Clear()

$ws.Range('ZZ16:ZZ17').NumberFormat = "@"
$ws.Range('ZZ16').Value = '13.93'
$ws.Range('ZZ17').Value = '20.43'
$ws.Range('ZZ16:ZZ17').Copy()
$ws.Range('D16:D17').PasteSpecial(-4104, $null, $null, $null)
$ws.Range('ZZ16:ZZ17').Clear()

$ws.Range('ZZ21:ZZ37').NumberFormat = "@"
$ws.Range('ZZ21').Value = '425.50'
$ws.Range('ZZ22').Value = '97.14'
$ws.Range('ZZ23').Value = '3.46'
$ws.Range('ZZ24').Value = '4.16'
$ws.Range('ZZ25').Value = '14.29'
$ws.Range('ZZ26').Value = '11.13'
$ws.Range('ZZ27').Value = '10.63'
$ws.Range('ZZ28').Value = '5.84'
$ws.Range('ZZ29').Value = '36.45'
$ws.Range('ZZ30').Value = '3.59'
$ws.Range('ZZ31').Value = '13.29'
$ws.Range('ZZ32').Value = '0.128'
$ws.Range('ZZ33').Value = '675.57'
$ws.Range('ZZ34').Value = '6.95'
$ws.Range('ZZ35').Value = '65.30'
$ws.Range('ZZ36').Value = '43.66'
$ws.Range('ZZ37').Value = '0.431'
$ws.Range('ZZ21:ZZ37').Copy()
$ws.Range('D21:D37').PasteSpecial(-4104, $null, $null, $null)
$ws.Range('ZZ21:ZZ37').Clear()

$ws.Range('ZZ40').NumberFormat = "@"
$ws.Range('ZZ40').Value = '3.40'
$ws.Range('ZZ40').Copy()
$ws.Range('D40').PasteSpecial(-4104, $null, $null, $null)
$ws.Range('ZZ40').Clear()

$ws.Range('ZZ43:ZZ51').NumberFormat = "@"
$ws.Range('ZZ43').Value = '0.0484'
$ws.Range('ZZ44').Value = '3.14'
$ws.Range('ZZ45').Value = '0.149'
$ws.Range('ZZ46').Value = '3.41'
$ws.Range('ZZ47').Value = '2.61'
$ws.Range('ZZ48').Value = '9.55'
$ws.Range('ZZ49').Value = '2.97'
$ws.Range('ZZ50').Value = '0.000270'
$ws.Range('ZZ51').Value = '145.48'
$ws.Range('ZZ43:ZZ51').Copy()
$ws.Range('D43:D51').PasteSpecial(-4104, $null, $null, $null)
$ws.Range('ZZ43:ZZ51').Clear()

$excel.CutCopyMode = $false
